$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" -------------------
# The shared string "Ready for handoff" is referenced by every Status-type
# cell across the three sheets (Overview!E/F rows 2-3, and the zh-cn / de-de
# sheets' Status column C rows 2-3). Update every occurrence so the text
# reads "In Translation" everywhere.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column widths: status columns shrink now that the text is shorter ------
# The report generator auto-sizes the Status columns to fit their text;
# now that the text reads "In Translation" instead of the longer
# "Ready for handoff", those columns narrow from ~17.22 chars to ~13.41 chars.
# "Overview" sheet: columns E (zh-cn) and F (de-de) hold the status text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# "zh-cn" / "de-de" sheets: column C holds the Status text.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
